# Daily refresh of the cryptos price/volume table (GitHub Actions bot).
#
# Columns: A=rank(0-based) B=Coin C=Link D=Price E=Volume(1h)
# All of B:E are stored as plain text in this sheet (prices use "."
# as a thousands separator in some rows, e.g. "35.214.03", and the
# volume column keeps its padded "  +1.31%  " formatting), so every
# write below goes through Set-TextCell, which forces the cell to
# Text before the assignment and then restores the "Normal" style so
# the workbook's styling is left untouched (Excel would otherwise
# silently re-interpret single-decimal strings like "239.59" as
# numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "35.214.03"
Set-TextCell "E2" "  +1.31%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.859.60"
Set-TextCell "E3" "  +1.67%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.40%  "

# Row 5 - BNB
Set-TextCell "D5" "239.59"
Set-TextCell "E5" "  +3.89%  "

# Row 6 - XRP
Set-TextCell "D6" "0.623"
Set-TextCell "E6" "  +0.90%  "

# Row 7 - USDC
Set-TextCell "E7" "  +0.40%  "

# Row 8 - Solana
Set-TextCell "D8" "42.28"
Set-TextCell "E8" "  +7.41%  "

# Row 9 - Cardano
Set-TextCell "E9" "  +0.92%  "

# Row 10 - Dogecoin
Set-TextCell "E10" "  +1.43%  "

# Row 11 - TRON
Set-TextCell "E11" "  +0.11%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell "D12" "2.128.21"
Set-TextCell "E12" "  +1.66%  "

# Row 13 - Chainlink
Set-TextCell "E13" "  +1.68%  "

# Row 14 - WrappedEther
Set-TextCell "D14" "1.858.84"
Set-TextCell "E14" "  +1.24%  "

# Row 15 - Polygon
Set-TextCell "D15" "0.677"
Set-TextCell "E15" "  +1.21%  "

# Row 16 - Polkadot
Set-TextCell "E16" "  +2.05%  "

# Row 17 - WrappedBTC
Set-TextCell "D17" "35.185.60"
Set-TextCell "E17" "  +1.17%  "

# Row 18 - Litecoin
Set-TextCell "D18" "69.91"
Set-TextCell "E18" "  +0.65%  "

# Row 19 - ShibaInu
Set-TextCell "E19" "  +1.46%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "240.60"

# Row 21 - Avalanche
Set-TextCell "D21" "12.23"
Set-TextCell "E21" "  +0.73%  "

# Row 22 - Uniswap
Set-TextCell "D22" "4.76"
Set-TextCell "E22" "  +1.93%  "

# Row 23 - Dai
Set-TextCell "E23" "  +0.35%  "

# Row 24 - Toncoin
Set-TextCell "E24" "  +0.86%  "

# Row 25 - Monero
Set-TextCell "D25" "169.48"
Set-TextCell "E25" "  -1.53%  "

# Row 26 - PancakeSwap
Set-TextCell "E26" "  +26.71%  "

# Row 27 - Cosmos
Set-TextCell "D27" "8.02"
Set-TextCell "E27" "  +3.51%  "

# Row 28 - EthereumClassic
Set-TextCell "D28" "17.66"
Set-TextCell "E28" "  +1.92%  "

# Row 29 - Stellar
Set-TextCell "E29" "  +0.25%  "

# Row 30 - BinanceUSD
Set-TextCell "E30" "  +0.34%  "

# Row 32 - Filecoin
Set-TextCell "E32" "  +2.04%  "

# Row 33 - WEMIXToken
Set-TextCell "E33" "  +27.51%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextCell "E34" "  +2.32%  "

# Row 35 - LidoDAOToken
Set-TextCell "E35" "  +10.32%  "

# Row 36 - ImmutableX
Set-TextCell "D36" "0.815"
Set-TextCell "E36" "  +16.86%  "

# Row 37 - TrustWalletToken
Set-TextCell "E37" "  +7.55%  "

# Row 38 - ARBITRUM
Set-TextCell "E38" "  +4.62%  "

# Row 39 - VeChain
Set-TextCell "E39" "  +4.12%  "

# Row 40 - Aave
Set-TextCell "D40" "89.94"
Set-TextCell "E40" "  -1.49%  "

# Row 41 - Maker
Set-TextCell "D41" "1.347.58"
Set-TextCell "E41" "  +0.64%  "

# Rows 42/43 - Kaspa and InjectiveProtocol swap ranking positions
Set-TextCell "B42" "Kaspa"
Set-TextCell "C42" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D42" "0.0596"
Set-TextCell "E42" "  +14.57%  "

Set-TextCell "B43" "InjectiveProtocol"
Set-TextCell "C43" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D43" "14.99"
Set-TextCell "E43" "  +3.52%  "

# Row 44 - RenderToken
Set-TextCell "D44" "2.33"
Set-TextCell "E44" "  +3.48%  "

# Row 45 - HuobiToken
Set-TextCell "D45" "2.43"
Set-TextCell "E45" "  +0.31%  "

# Row 46 - Gas
Set-TextCell "D46" "12.41"
Set-TextCell "E46" "  +42.42%  "

# Row 47 - MXToken
Set-TextCell "E47" "  -0.53%  "

# Row 48 - FraxShare
Set-TextCell "D48" "6.60"
Set-TextCell "E48" "  +5.59%  "

# Row 49 - RocketPoolETH
Set-TextCell "D49" "2.044.78"
Set-TextCell "E49" "  +1.81%  "

# Row 50 - Cronos
Set-TextCell "D50" "0.0681"
Set-TextCell "E50" "  +1.42%  "

# Row 51 - PaxDollar
Set-TextCell "E51" "  +0.44%  "
